# Commit: "Using median instead of mean"
#
# The summary rows (36, 73, 110) for each of the four column-groups
# (C:F, J:M, Q:T, X:AA) used AVERAGE() over the block of raw samples above
# them. Switch those formulas to MEDIAN() over the same ranges, leaving the
# AVEDEV() rows (37, 74, 111) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("C", "D", "E", "F", "J", "K", "L", "M", "Q", "R", "S", "T", "X", "Y", "Z", "AA")

# (summary row, first data row, last data row)
$blocks = @(
    @(36, 6, 35),
    @(73, 43, 72),
    @(110, 80, 109)
)

foreach ($block in $blocks) {
    $summaryRow = $block[0]
    $firstRow = $block[1]
    $lastRow = $block[2]

    foreach ($col in $columns) {
        $cell = $ws.Range("$col$summaryRow")
        $cell.Formula = "=MEDIAN($col$firstRow`:$col$lastRow)"
    }
}
